$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated quote values (Dólar, Euro, Bitcoin)
$ws.Range("B2").Value = 5.8576
$ws.Range("B3").Value = 6.65779
$ws.Range("B4").Value = 502853000

# Updated "last update" timestamp for all three rows
$newDate = 45759.91172538449
$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
